$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = [double]"0"
$ws.Range("F2").Value = [double]"24.09000000000033"
$ws.Range("H2").Value = [double]"0.02449221832896331"
$ws.Range("I2").Value = [double]"0.02449221832896331"
$ws.Range("L2").Value = [double]"35.70269648742626"
$ws.Range("M2").Value = "[2.499768982014473, 68.90562399283804]"
$ws.Range("N2").Value = [double]"0.03566902510678172"
$ws.Range("O2").Value = [double]"0.03566902510678172"
$ws.Range("P2").Value = [double]"1.113237036407194"
$ws.Range("Q2").Value = "[-0.11950002085726741, 2.345974093671656]"
$ws.Range("R2").Value = [double]"0.07559135083873225"
$ws.Range("S2").Value = [double]"0.07559135083873225"
$ws.Range("T2").Value = [double]"69.44369302820272"
$ws.Range("U2").Value = "[51.54332208470318, 87.34406397170225]"
$ws.Range("V2").Value = [double]"6.446092548628712e-10"
$ws.Range("W2").Value = [double]"6.446092548628712e-10"
$ws.Range("X2").Value = [double]"19.82180180180207"
$ws.Range("Y2").Value = [double]"15.09543543543564"
$ws.Range("Z2").Value = [double]"24.54816816816849"
# Row 3
$ws.Range("B3").Value = [double]"0"
$ws.Range("F3").Value = [double]"24.09000000000033"
$ws.Range("H3").Value = [double]"0.03986696404325374"
$ws.Range("I3").Value = [double]"0.03986696404325374"
$ws.Range("L3").Value = [double]"26.28915878074509"
$ws.Range("M3").Value = "[-1.4402092107369597, 54.01852677222714]"
$ws.Range("N3").Value = [double]"0.06258609044086905"
$ws.Range("O3").Value = [double]"0.06258609044086905"
$ws.Range("P3").Value = [double]"1.088079137279347"
$ws.Range("Q3").Value = "[-0.45913165908319353, 2.6352899336418885]"
$ws.Range("R3").Value = [double]"0.1635373886001132"
$ws.Range("S3").Value = [double]"0.1635373886001132"
$ws.Range("T3").Value = [double]"49.28142965199198"
$ws.Range("U3").Value = "[35.06780757907862, 63.49505172490535]"
$ws.Range("V3").Value = [double]"1.078808553600652e-08"
$ws.Range("W3").Value = [double]"1.078808553600652e-08"
$ws.Range("X3").Value = [double]"19.91825825825853"
$ws.Range("Y3").Value = [double]"13.98618618618638"
$ws.Range("Z3").Value = [double]"25.85033033033068"
# Row 4
$ws.Range("F4").Value = [double]"24.09000000000033"
$ws.Range("H4").Value = [double]"1.51609515663953e-06"
$ws.Range("I4").Value = [double]"1.51609515663953e-06"
$ws.Range("L4").Value = [double]"60.52603728802253"
$ws.Range("M4").Value = "[33.61713634233118, 87.43493823371388]"
$ws.Range("N4").Value = [double]"4.314437544272742e-05"
$ws.Range("O4").Value = [double]"4.314437544272742e-05"
$ws.Range("P4").Value = [double]"1.062921238151501"
$ws.Range("Q4").Value = "[0.6100790538502707, 1.515763422452732]"
$ws.Range("R4").Value = [double]"2.265579905924575e-05"
$ws.Range("S4").Value = [double]"2.265579905924575e-05"
$ws.Range("T4").Value = [double]"63.49941184595193"
$ws.Range("U4").Value = "[49.248045100746936, 77.75077859115692]"
$ws.Range("V4").Value = [double]"1.378119840467207e-11"
$ws.Range("W4").Value = [double]"1.378119840467207e-11"
$ws.Range("X4").Value = [double]"20.01471471471499"
$ws.Range("Y4").Value = [double]"18.27849849849875"
$ws.Range("Z4").Value = [double]"21.75093093093123"
# Row 5
$ws.Range("F5").Value = [double]"24.09000000000033"
$ws.Range("H5").Value = [double]"5.781573336405277e-05"
$ws.Range("I5").Value = [double]"5.781573336405277e-05"
$ws.Range("L5").Value = [double]"55.108209121087"
$ws.Range("M5").Value = "[24.29301113105305, 85.92340711112095]"
$ws.Range("N5").Value = [double]"0.0007849498755945206"
$ws.Range("O5").Value = [double]"0.0007849498755945206"
$ws.Range("P5").Value = [double]"1.075500187715424"
$ws.Range("Q5").Value = "[0.5094474573388847, 1.6415529180919641]"
$ws.Range("R5").Value = [double]"0.0003980332039408285"
$ws.Range("S5").Value = [double]"0.0003980332039408285"
$ws.Range("T5").Value = [double]"65.67372477419154"
$ws.Range("U5").Value = "[49.78502086875025, 81.56242867963283]"
$ws.Range("V5").Value = [double]"1.165463281438406e-10"
$ws.Range("W5").Value = [double]"1.165463281438406e-10"
$ws.Range("X5").Value = [double]"19.96648648648676"
$ws.Range("Y5").Value = [double]"17.79621621621646"
$ws.Range("Z5").Value = [double]"22.13675675675706"
# Row 6
$ws.Range("B6").Value = [double]"1"
$ws.Range("F6").Value = [double]"25.64000000000057"
$ws.Range("H6").Value = [double]"0.006372194202847847"
$ws.Range("I6").Value = [double]"0.006372194202847847"
$ws.Range("L6").Value = [double]"40.12383009085594"
$ws.Range("M6").Value = "[7.945656040678969, 72.30200414103291]"
$ws.Range("N6").Value = [double]"0.01568186423392004"
$ws.Range("O6").Value = [double]"0.01568186423392004"
$ws.Range("P6").Value = [double]"1.000026490331886"
$ws.Range("Q6").Value = "[0.20755266780473125, 1.7925003128590413]"
$ws.Range("R6").Value = [double]"0.01454917398794375"
$ws.Range("S6").Value = [double]"0.01454917398794375"
$ws.Range("T6").Value = [double]"56.51221903102515"
$ws.Range("U6").Value = "[39.55777934247877, 73.46665871957154]"
$ws.Range("V6").Value = [double]"2.714983859952724e-08"
$ws.Range("W6").Value = [double]"2.714983859952724e-08"
$ws.Range("X6").Value = [double]"21.55915915915964"
$ws.Range("Y6").Value = [double]"18.32528528528569"
$ws.Range("Z6").Value = [double]"24.79303303303358"
# Row 7
$ws.Range("F7").Value = [double]"25.64000000000057"
$ws.Range("H7").Value = [double]"3.166039885837968e-05"
$ws.Range("I7").Value = [double]"3.166039885837968e-05"
$ws.Range("L7").Value = [double]"68.98050167449298"
$ws.Range("M7").Value = "[38.8740825111926, 99.08692083779337]"
$ws.Range("N7").Value = [double]"3.277639638055163e-05"
$ws.Range("O7").Value = [double]"3.277639638055163e-05"
$ws.Range("P7").Value = [double]"0.04402632347373014"
$ws.Range("Q7").Value = "[-0.48428955821104047, 0.5723422051585008]"
$ws.Range("R7").Value = [double]"0.8674596780438755"
$ws.Range("S7").Value = [double]"0.8674596780438755"
$ws.Range("T7").Value = [double]"79.51357009938269"
$ws.Range("U7").Value = "[60.970038993552166, 98.05710120521321]"
$ws.Range("V7").Value = [double]"4.163291933423352e-11"
$ws.Range("W7").Value = [double]"4.163291933423352e-11"
$ws.Range("X7").Value = [double]"25.46034034034091"
$ws.Range("Y7").Value = [double]"23.30442442442494"
$ws.Range("Z7").Value = [double]"27.61625625625688"
# Row 8
$ws.Range("F8").Value = [double]"25.64000000000057"
$ws.Range("H8").Value = [double]"0.001197459377738741"
$ws.Range("I8").Value = [double]"0.001197459377738741"
$ws.Range("L8").Value = [double]"42.6496446289797"
$ws.Range("M8").Value = "[17.885594287673356, 67.41369497028604]"
$ws.Range("N8").Value = [double]"0.001163413547255354"
$ws.Range("O8").Value = [double]"0.001163413547255354"
$ws.Range("P8").Value = [double]"-0.0503157982556921"
$ws.Range("Q8").Value = "[-0.792473822527155, 0.6918422260157708]"
$ws.Range("R8").Value = [double]"0.8919960697227634"
$ws.Range("S8").Value = [double]"0.8919960697227634"
$ws.Range("T8").Value = [double]"54.13073879856472"
$ws.Range("U8").Value = "[39.20658092187399, 69.05489667525545]"
$ws.Range("V8").Value = [double]"3.600878040188604e-09"
$ws.Range("W8").Value = [double]"3.600878040188604e-09"
$ws.Range("X8").Value = [double]"0.205325325325326"
$ws.Range("Y8").Value = [double]"-2.823223223223291"
$ws.Range("Z8").Value = [double]"3.233873873873943"
# Row 9
$ws.Range("F9").Value = [double]"25.64000000000057"
$ws.Range("H9").Value = [double]"0.00034899771868524"
$ws.Range("I9").Value = [double]"0.00034899771868524"
$ws.Range("L9").Value = [double]"45.42906129931608"
$ws.Range("M9").Value = "[15.989439552448133, 74.86868304618403]"
$ws.Range("N9").Value = [double]"0.003260163494942292"
$ws.Range("O9").Value = [double]"0.003260163494942292"
$ws.Range("P9").Value = [double]"-0.9685791164220783"
$ws.Range("Q9").Value = "[-1.58494764505431, -0.3522105877898465]"
$ws.Range("R9").Value = [double]"0.002780585060515151"
$ws.Range("S9").Value = [double]"0.002780585060515151"
$ws.Range("T9").Value = [double]"65.16988928371256"
$ws.Range("U9").Value = "[49.9362609797662, 80.40351758765892]"
$ws.Range("V9").Value = [double]"4.446110146716364e-11"
$ws.Range("W9").Value = [double]"4.446110146716364e-11"
$ws.Range("X9").Value = [double]"3.952512512512602"
$ws.Range("Y9").Value = [double]"1.43727727727731"
$ws.Range("Z9").Value = [double]"6.467747747747893"
# Row 10
$ws.Range("F10").Value = [double]"25.64000000000057"
$ws.Range("H10").Value = [double]"0.0001699394381059482"
$ws.Range("I10").Value = [double]"0.0001699394381059482"
$ws.Range("L10").Value = [double]"48.4936482495772"
$ws.Range("M10").Value = "[19.07547689672042, 77.91181960243398]"
$ws.Range("N10").Value = [double]"0.001790694547753935"
$ws.Range("O10").Value = [double]"0.001790694547753935"
$ws.Range("P10").Value = [double]"-0.6792632764518469"
$ws.Range("Q10").Value = "[-1.2578949563923096, -0.1006315965113842]"
$ws.Range("R10").Value = [double]"0.02243492632858546"
$ws.Range("S10").Value = [double]"0.02243492632858546"
$ws.Range("T10").Value = [double]"57.25046358743437"
$ws.Range("U10").Value = "[42.04267153308679, 72.45825564178196]"
$ws.Range("V10").Value = [double]"1.407627348015694e-09"
$ws.Range("W10").Value = [double]"1.407627348015694e-09"
$ws.Range("X10").Value = [double]"2.771891891891954"
$ws.Range("Y10").Value = [double]"0.4106506506506591"
$ws.Range("Z10").Value = [double]"5.133133133133249"
